$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 43-64)
$ws.Range("A43:B64").EntireRow.Delete() | Out-Null

# Refreshed dataset: Date (col A) / CombinedVaderSentiment (col B)
$dates = @(
    44543,
    44550,
    44557,
    44564,
    44571,
    44578,
    44585,
    44592,
    44599,
    44606,
    44613,
    44620,
    44627,
    44634,
    44641,
    44648,
    44655,
    44662,
    44669,
    44676,
    44683,
    44690,
    44697,
    44704,
    44711,
    44718,
    44725,
    44732,
    44739,
    44746,
    44753,
    44760,
    44767,
    44774,
    44781,
    44788,
    44795,
    44802,
    44809,
    44816,
    44823
)
$values = @(
    1.3612,
    0.56335,
    0.5781000000000001,
    0.852,
    0.995,
    0.9924999999999999,
    0.995,
    0.995,
    0.986,
    0.8704499999999999,
    0.9975000000000001,
    1.1598,
    0.998,
    0.19515,
    0.86225,
    0.957,
    0.998,
    0.967,
    0.9804999999999999,
    0.9895,
    1.0248,
    0.679,
    0.999,
    0.86465,
    0.996,
    0.956,
    0.957,
    0.988,
    0.993,
    0.98,
    0.998,
    0.997,
    0.995,
    0.90655,
    0.7800499999999999,
    0.9955000000000001,
    0.99,
    0.9924999999999999,
    0.986,
    0.993,
    0.998
)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

Write-Host "Updated $($dates.Count) rows"
